$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "DCase 2018 Baseline" (TUT 2017) result block, rows 3-6, columns C:F
$ws.Range("C3").Value = "0.8763"
$ws.Range("D3").Value = "0.8879"
$ws.Range("E3").Value = "0.8622"
$ws.Range("F3").Value = "0.8011"

$ws.Range("C4").Value = "0.4573"
$ws.Range("D4").Value = "0.6865"
$ws.Range("E4").Value = "0.6313"
$ws.Range("F4").Value = "0.5528"

$ws.Range("C5").Value = "0.3894"
$ws.Range("D5").Value = "0.5071"
$ws.Range("E5").Value = "0.5030"
$ws.Range("F5").Value = "0.3630"

$ws.Range("C6").Value = "0.4206"
$ws.Range("D6").Value = "0.5833"
$ws.Range("E6").Value = "0.5030"
$ws.Range("F6").Value = "0.4382"

# Row 6 (F1-Score for DCase baseline) becomes bold
$ws.Range("C6:F6").Font.Bold = $true

# Row 11 (F1-Score for AdvancedCRNN) is no longer bold
$ws.Range("C11:F11").Font.Bold = $false

# Move selection to F11
$ws.Range("F11").Select()
